$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The registration algorithm no longer collects a date of birth, so the
# "dateOfBirth" input-type row is cleared out...
$ws.Range("B8").ClearContents()

# ...and the final "save" step is updated to drop dateOfBirth from the
# list of fields persisted to the database.
$ws.Range("B30").Value = "save name, email, userName, password to datebase"

# Leave the view scrolled to / focused on the line that was just edited.
$excel.ActiveWindow.ScrollRow = 21
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B30").Select()
